$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "HEALTHY REWARDS`nOFFER WITH CARD`n`$1499/lb. 2`nSAVE `$10/lb.`nFresh, Wild`nNON-GMO`nFRESH`nSockeye`nSalmon Fillets`n100% Traceable and Sustainably Sourced, All Earth Fare`nSalmon is non-GMO, Superior Flavor, Great Source of`nOmega-3 Fatty Acids, Product of the USA (Alaska)`nDiscount Taken at Register"

$ws.Range("A2").Value = $newText

$ws.Rows.Item(2).RowHeight = 187.2
